$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows (6,7,8) to make room for new log entries ---
$ws.Rows("6:8").Insert()

# --- Copy formatting baseline (borders etc.) from row 4 into new rows 6-8 ---
$ws.Range("A4:F4").Copy()
$ws.Range("A6:F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update row 4: "Setup the Project" ---
$ws.Range("A4").Value = "Setup the Project "
$ws.Range("B4").Value = 0.25
$ws.Range("C4").Value = 45055
$ws.Range("D4").Value = "Setup the github repo and create an empty unity 3d Project"

# --- Update row 5: "Feature 1: Setup Grid generation" ---
$ws.Range("A5").Value = "Feature 1: Setup Grid generation"
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = 45055
$ws.Range("D5").Value = "* Create a button to trigger generation `n* Create a Base prefab (Puzzle unit), in prefab set offset to the cube of 0.5 `nunit on Oy axis so that I can easily position cubes on bottom edge instead of center`n* Create a primitive script setup of 2D byte array representing the metadata for the maze ( 1 -> create cube, 0 -> will skip the cube) `n* Created a simple function that would create the cubes based on array values ( for now just treat like it has only 1s so have a plane of cubes). This would allow me to visually debug the next steps.`n"

# --- Update row 6: "Feature 1: Create the first iteration of the Maze" ---
$ws.Range("A6").Value = "Feature 1: Create the first iteration of the Maze"
$ws.Range("B6").Value = 2
$ws.Cells.Item(6,3).Value2 = "14/5/2023"
$ws.Range("D6").Value = "* Using the Randomized Prim's pseudo algorith from the provided wiki page created the first geneartion. `n* For now did not bother with code quality, focused on the core generation.`n"

# --- D5:D8 need wrap text (long descriptions / placeholder rows) ---
$ws.Range("D5:D8").WrapText = $true

# --- Number format: Amount of hours column now shows 2 decimals ---
$ws.Range("B4:B31").NumberFormat = "0.00"
$ws.Columns("B").NumberFormat = "0.00"

# --- Row heights ---
$ws.Rows(2).RowHeight = 57.6
$ws.Rows(5).RowHeight = 93
$ws.Rows(6).RowHeight = 49.2
$ws.Rows(7).RowHeight = 106.2
$ws.Rows(8).RowHeight = 17.4

# --- Column width tweaks ---
$ws.Columns("A").ColumnWidth = 24.765625
$ws.Columns("B").ColumnWidth = 10.765625
$ws.Columns("C").ColumnWidth = 8.61328125
$ws.Columns("D").ColumnWidth = 38.61328125

# --- View: selection / scroll ---
$ws.Range("D7").Select()
try { $excel.ActiveWindow.ScrollRow = 2 } catch {}

Write-Host "done"